$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: E2 (zh-cn status) and F2 (de-de status)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: C2 (Status column)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: C2 (Status column)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the status columns to fit the new, shorter text ---
# ColumnWidth is expressed in character units; the saved column width is
# quantized to the nearest 1/6th of a character unit by the engine, so pick
# the input that lands on the closest achievable width to 13.4101845877511.
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
